# Append new shipping/pipelay log entries (rows 67-131) to Sheet1.
# Mirrors an upload that extended the activity log from 2024-06-01
# through 2024-06-20 (dates stored as serials 45445-45463).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 65,6
$data[0,0] = 45445
$data[0,1] = "00:00:00"
$data[0,2] = "02:27:00"
$data[0,3] = "10_DAPLD_6km"
$data[0,4] = "PL"
$data[0,5] = "Laying"
$data[1,0] = 45445
$data[1,1] = "02:27:00"
$data[1,2] = "10:21:00"
$data[1,3] = "10_DAPLD_6km"
$data[1,4] = "PL"
$data[1,5] = "Buckle detector"
$data[2,0] = 45445
$data[2,1] = "10:21:00"
$data[2,2] = "23:59:00"
$data[2,3] = "10_DAPLD_6km"
$data[2,4] = "PL"
$data[2,5] = "Laying"
$data[3,0] = 45446
$data[3,1] = "00:00:00"
$data[3,2] = "23:59:00"
$data[3,3] = "10_DAPLD_6km"
$data[3,4] = "PL"
$data[3,5] = "Laying"
$data[4,0] = 45447
$data[4,1] = "00:00:00"
$data[4,2] = "20:05:00"
$data[4,3] = "10_DAPLD_6km"
$data[4,4] = "PL"
$data[4,5] = "Laying"
$data[5,0] = 45447
$data[5,1] = "20:05:00"
$data[5,2] = "23:59:00"
$data[5,3] = "10_DAPLD_6km"
$data[5,4] = "PL"
$data[5,5] = "Lay-down"
$data[6,0] = 45448
$data[6,1] = "00:00:00"
$data[6,2] = "07:22:00"
$data[6,3] = "10_DAPLD_6km"
$data[6,4] = "PL"
$data[6,5] = "Lay-down"
$data[7,0] = 45448
$data[7,1] = "07:22:00"
$data[7,2] = "10:50:00"
$data[7,3] = "10_DAPLD_6km"
$data[7,4] = "Others"
$data[7,5] = "-"
$data[8,0] = 45448
$data[8,1] = "10:50:00"
$data[8,2] = "17:53:00"
$data[8,3] = "10_DAPLD_6km"
$data[8,4] = "Anchor"
$data[8,5] = "-"
$data[9,0] = 45448
$data[9,1] = "17:53:00"
$data[9,2] = "23:59:00"
$data[9,3] = "10_DAPLD_6km"
$data[9,4] = "Downtime"
$data[9,5] = "AHT/Barge"
$data[10,0] = 45449
$data[10,1] = "00:00:00"
$data[10,2] = "01:35:00"
$data[10,3] = "10_DAPLD_6km"
$data[10,4] = "Downtime"
$data[10,5] = "AHT/Barge"
$data[11,0] = 45449
$data[11,1] = "01:35:00"
$data[11,2] = "04:30:00"
$data[11,3] = "10_DAPLD_6km"
$data[11,4] = "Anchor"
$data[11,5] = "-"
$data[12,0] = 45449
$data[12,1] = "04:31:00"
$data[12,2] = "04:33:00"
$data[12,3] = "10_DAPLD_6km"
$data[12,4] = "Anchor"
$data[12,5] = "-"
$data[13,0] = 45449
$data[13,1] = "04:33:00"
$data[13,2] = "05:35:00"
$data[13,3] = "10_DAPLD_6km"
$data[13,4] = "Anchor"
$data[13,5] = "-"
$data[14,0] = 45449
$data[14,1] = "05:35:00"
$data[14,2] = "07:35:00"
$data[14,3] = "10_DAPLD_6km"
$data[14,4] = "Downtime"
$data[14,5] = "AHT/Barge"
$data[15,0] = 45449
$data[15,1] = "07:35:00"
$data[15,2] = "09:10:00"
$data[15,3] = "10_DAPLD_6km"
$data[15,4] = "Anchor"
$data[15,5] = "-"
$data[16,0] = 45449
$data[16,1] = "09:10:00"
$data[16,2] = "18:40:00"
$data[16,3] = "Transit"
$data[16,4] = "Field_Move"
$data[16,5] = "-"
$data[17,0] = 45449
$data[17,1] = "18:40:00"
$data[17,2] = "23:59:00"
$data[17,3] = "10_GOPLE_7km"
$data[17,4] = "Anchor"
$data[17,5] = "-"
$data[18,0] = 45450
$data[18,1] = "00:00:00"
$data[18,2] = "03:30:00"
$data[18,3] = "10_GOPLE_7km"
$data[18,4] = "Anchor"
$data[18,5] = "-"
$data[19,0] = 45450
$data[19,1] = "03:30:00"
$data[19,2] = "21:34:00"
$data[19,3] = "10_GOPLE_7km"
$data[19,4] = "PL"
$data[19,5] = "Start-up"
$data[20,0] = 45450
$data[20,1] = "21:34:00"
$data[20,2] = "23:59:00"
$data[20,3] = "10_GOPLE_7km"
$data[20,4] = "PL"
$data[20,5] = "Laying"
$data[21,0] = 45451
$data[21,1] = "00:00:00"
$data[21,2] = "03:50:00"
$data[21,3] = "10_GOPLE_7km"
$data[21,4] = "PL"
$data[21,5] = "Laying"
$data[22,0] = 45451
$data[22,1] = "03:50:00"
$data[22,2] = "05:38:00"
$data[22,3] = "10_GOPLE_7km"
$data[22,4] = "PL"
$data[22,5] = "Buckle detector"
$data[23,0] = 45451
$data[23,1] = "05:38:00"
$data[23,2] = "23:59:00"
$data[23,3] = "10_GOPLE_7km"
$data[23,4] = "PL"
$data[23,5] = "Laying"
$data[24,0] = 45452
$data[24,1] = "00:00:00"
$data[24,2] = "08:07:00"
$data[24,3] = "10_GOPLE_7km"
$data[24,4] = "PL"
$data[24,5] = "Laying"
$data[25,0] = 45452
$data[25,1] = "08:07:00"
$data[25,2] = "08:31:00"
$data[25,3] = "10_GOPLE_7km"
$data[25,4] = "Anchor"
$data[25,5] = "-"
$data[26,0] = 45452
$data[26,1] = "08:31:00"
$data[26,2] = "12:25:00"
$data[26,3] = "10_GOPLE_7km"
$data[26,4] = "PL"
$data[26,5] = "Laying"
$data[27,0] = 45452
$data[27,1] = "12:25:00"
$data[27,2] = "13:42:00"
$data[27,3] = "10_GOPLE_7km"
$data[27,4] = "Anchor"
$data[27,5] = "-"
$data[28,0] = 45452
$data[28,1] = "13:42:00"
$data[28,2] = "21:24:00"
$data[28,3] = "10_GOPLE_7km"
$data[28,4] = "PL"
$data[28,5] = "Laying"
$data[29,0] = 45452
$data[29,1] = "21:24:00"
$data[29,2] = "21:33:00"
$data[29,3] = "10_GOPLE_7km"
$data[29,4] = "Anchor"
$data[29,5] = "-"
$data[30,0] = 45452
$data[30,1] = "21:33:00"
$data[30,2] = "23:59:00"
$data[30,3] = "10_GOPLE_7km"
$data[30,4] = "PL"
$data[30,5] = "Laying"
$data[31,0] = 45453
$data[31,1] = "00:00:00"
$data[31,2] = "01:00:00"
$data[31,3] = "10_GOPLE_7km"
$data[31,4] = "PL"
$data[31,5] = "Laying"
$data[32,0] = 45453
$data[32,1] = "01:00:00"
$data[32,2] = "07:15:00"
$data[32,3] = "10_GOPLE_7km"
$data[32,4] = "WOW"
$data[32,5] = "Hs >2m"
$data[33,0] = 45453
$data[33,1] = "07:15:00"
$data[33,2] = "23:59:00"
$data[33,3] = "10_GOPLE_7km"
$data[33,4] = "PL"
$data[33,5] = "Laying"
$data[34,0] = 45454
$data[34,1] = "00:00:00"
$data[34,2] = "22:13:00"
$data[34,3] = "10_GOPLE_7km"
$data[34,4] = "PL"
$data[34,5] = "Laying"
$data[35,0] = 45454
$data[35,1] = "22:13:00"
$data[35,2] = "23:59:00"
$data[35,3] = "10_GOPLE_7km"
$data[35,4] = "Downtime"
$data[35,5] = "AHT/Barge"
$data[36,0] = 45455
$data[36,1] = "00:00:00"
$data[36,2] = "09:27:00"
$data[36,3] = "10_GOPLE_7km"
$data[36,4] = "Downtime"
$data[36,5] = "AHT/Barge"
$data[37,0] = 45455
$data[37,1] = "09:27:00"
$data[37,2] = "09:45:00"
$data[37,3] = "10_GOPLE_7km"
$data[37,4] = "PL"
$data[37,5] = "Laying"
$data[38,0] = 45455
$data[38,1] = "09:45:00"
$data[38,2] = "23:59:00"
$data[38,3] = "10_GOPLE_7km"
$data[38,4] = "Downtime"
$data[38,5] = "AHT/Barge"
$data[39,0] = 45456
$data[39,1] = "00:00:00"
$data[39,2] = "08:01:00"
$data[39,3] = "10_GOPLE_7km"
$data[39,4] = "Downtime"
$data[39,5] = "AHT/Barge"
$data[40,0] = 45456
$data[40,1] = "08:01:00"
$data[40,2] = "08:18:00"
$data[40,3] = "10_GOPLE_7km"
$data[40,4] = "PL"
$data[40,5] = "Laying"
$data[41,0] = 45456
$data[41,1] = "08:18:00"
$data[41,2] = "17:04:00"
$data[41,3] = "10_GOPLE_7km"
$data[41,4] = "Downtime"
$data[41,5] = "AHT/Barge"
$data[42,0] = 45456
$data[42,1] = "17:04:00"
$data[42,2] = "23:59:00"
$data[42,3] = "10_GOPLE_7km"
$data[42,4] = "PL"
$data[42,5] = "Laying"
$data[43,0] = 45457
$data[43,1] = "00:00:00"
$data[43,2] = "06:07:00"
$data[43,3] = "10_GOPLE_7km"
$data[43,4] = "PL"
$data[43,5] = "Lay-down"
$data[44,0] = 45457
$data[44,1] = "06:07:00"
$data[44,2] = "06:47:00"
$data[44,3] = "10_GOPLE_7km"
$data[44,4] = "Anchor"
$data[44,5] = "-"
$data[45,0] = 45457
$data[45,1] = "06:47:00"
$data[45,2] = "14:00:00"
$data[45,3] = "10_GOPLE_7km"
$data[45,4] = "PL"
$data[45,5] = "Lay-down"
$data[46,0] = 45457
$data[46,1] = "14:00:00"
$data[46,2] = "16:41:00"
$data[46,3] = "10_GOPLE_7km"
$data[46,4] = "PL"
$data[46,5] = "Others"
$data[47,0] = 45457
$data[47,1] = "16:41:00"
$data[47,2] = "23:59:00"
$data[47,3] = "10_GOPLE_7km"
$data[47,4] = "Anchor"
$data[47,5] = "-"
$data[48,0] = 45458
$data[48,1] = "00:00:00"
$data[48,2] = "02:10:00"
$data[48,3] = "10_GOPLE_7km"
$data[48,4] = "Anchor"
$data[48,5] = "-"
$data[49,0] = 45458
$data[49,1] = "02:10:00"
$data[49,2] = "21:10:00"
$data[49,3] = "Transit"
$data[49,4] = "Field_Move"
$data[49,5] = "-"
$data[50,0] = 45458
$data[50,1] = "21:10:00"
$data[50,2] = "23:59:00"
$data[50,3] = "16_PLPLQ_3km"
$data[50,4] = "Anchor"
$data[50,5] = "-"
$data[51,0] = 45459
$data[51,1] = "00:00:00"
$data[51,2] = "08:48:00"
$data[51,3] = "16_PLPLQ_3km"
$data[51,4] = "Anchor"
$data[51,5] = "-"
$data[52,0] = 45459
$data[52,1] = "08:48:00"
$data[52,2] = "16:04:00"
$data[52,3] = "16_PLPLQ_3km"
$data[52,4] = "PL"
$data[52,5] = "Start-up"
$data[53,0] = 45459
$data[53,1] = "16:04:00"
$data[53,2] = "23:59:00"
$data[53,3] = "16_PLPLQ_3km"
$data[53,4] = "PL"
$data[53,5] = "Laying"
$data[54,0] = 45460
$data[54,1] = "00:00:00"
$data[54,2] = "11:37:00"
$data[54,3] = "16_PLPLQ_3km"
$data[54,4] = "PL"
$data[54,5] = "Laying"
$data[55,0] = 45460
$data[55,1] = "11:37:00"
$data[55,2] = "16:00:00"
$data[55,3] = "16_PLPLQ_3km"
$data[55,4] = "PL"
$data[55,5] = "Buckle detector"
$data[56,0] = 45460
$data[56,1] = "16:00:00"
$data[56,2] = "23:59:00"
$data[56,3] = "16_PLPLQ_3km"
$data[56,4] = "PL"
$data[56,5] = "Laying"
$data[57,0] = 45461
$data[57,1] = "00:00:00"
$data[57,2] = "21:10:00"
$data[57,3] = "16_PLPLQ_3km"
$data[57,4] = "PL"
$data[57,5] = "Laying"
$data[58,0] = 45461
$data[58,1] = "21:10:00"
$data[58,2] = "23:59:00"
$data[58,3] = "16_PLPLQ_3km"
$data[58,4] = "Anchor"
$data[58,5] = "-"
$data[59,0] = 45462
$data[59,1] = "00:00:00"
$data[59,2] = "02:03:00"
$data[59,3] = "16_PLPLQ_3km"
$data[59,4] = "Anchor"
$data[59,5] = "-"
$data[60,0] = 45462
$data[60,1] = "02:03:00"
$data[60,2] = "03:29:00"
$data[60,3] = "16_PLPLQ_3km"
$data[60,4] = "PL"
$data[60,5] = "Laying"
$data[61,0] = 45462
$data[61,1] = "03:29:00"
$data[61,2] = "15:30:00"
$data[61,3] = "16_PLPLQ_3km"
$data[61,4] = "PL"
$data[61,5] = "Lay-down"
$data[62,0] = 45462
$data[62,1] = "15:30:00"
$data[62,2] = "19:50:00"
$data[62,3] = "16_PLPLQ_3km"
$data[62,4] = "PL"
$data[62,5] = "Others"
$data[63,0] = 45462
$data[63,1] = "19:50:00"
$data[63,2] = "23:59:00"
$data[63,3] = "16_PLPLQ_3km"
$data[63,4] = "Anchor"
$data[63,5] = "-"
$data[64,0] = 45463
$data[64,1] = "00:00:00"
$data[64,2] = "05:00:00"
$data[64,3] = "16_PLPLQ_3km"
$data[64,4] = "Anchor"
$data[64,5] = "-"

$ws.Range("A67:F131").Value = $data

# Match the author's final view state: scrolled to row 106, E110 selected.
$ws.Range("E110").Select()
$excel.ActiveWindow.ScrollRow = 106
